{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paras.items[paras.items.length - 1];\n\nconst list = lastParagraph.list;\nlist.load(\"id\");\nawait context.sync();\nconst listId = list.id;\n\nconst newTexts = [\n  \"Now place sign out form in layout as we need the sign out button to be included on all pages. \",\n  \"Wrote a delete request in the controller. \",\n  \"Test still failing as it says it can\\u2019t find the sign out button. I\\u2019m racking up\",\n  \"Kept failing, couldn\\u2019t see the sign out button, so I created a let statement at the top of the feature allowing factory girl to build the user, and then a before each statement calling the sign up and sign in methods.\"\n];\n\nlet anchor = lastParagraph;\nfor (const t of newTexts) {\n  const p = anchor.insertParagraph(t, \"After\");\n  p.style = \"List Paragraph\";\n  p.attachToList(listId, 0);\n  anchor = p;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$newTexts = @(\n  \"Now place sign out form in layout as we need the sign out button to be included on all pages. \",\n  \"Wrote a delete request in the controller. \",\n  \"Test still failing as it says it can\u2019t find the sign out button. I\u2019m racking up\",\n  \"Kept failing, couldn\u2019t see the sign out button, so I created a let statement at the top of the feature allowing factory girl to build the user, and then a before each statement calling the sign up and sign in methods.\"\n)\n\nforeach ($t in $newTexts) {\n  $p = $d.Paragraphs.Last\n  $p.Range.InsertParagraphAfter()\n  $newP = $d.Paragraphs.Last\n  $newP.Range.InsertAfter($t)\n}\n"}
